$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ordered attribute/type rows (rows 3-14 and 16-20); rows 1,2,15,21 unchanged.
$ws.Range("A3").Value = "case:concept:name"
$ws.Range("B3").Value = "str"

$ws.Range("A4").Value = "org:resource"
$ws.Range("B4").Value = "str"

$ws.Range("A5").Value = "response_status_code"
$ws.Range("B5").Value = "float"

$ws.Range("A6").Value = "human_workstation_green_button_pressed"
$ws.Range("B6").Value = "float"

$ws.Range("A7").Value = "identifier:id"
$ws.Range("B7").Value = "str"

$ws.Range("A8").Value = "lifecycle:state"
$ws.Range("B8").Value = "str"

$ws.Range("A9").Value = "unsatisfied_condition_description"
$ws.Range("B9").Value = "str"

$ws.Range("A10").Value = "parameters"
$ws.Range("B10").Value = "dict"

$ws.Range("A11").Value = "event_id"
$ws.Range("B11").Value = "str"

$ws.Range("A12").Value = "planned_operation_time"
$ws.Range("B12").Value = "str"

$ws.Range("A13").Value = "complete_service_time"
$ws.Range("B13").Value = "str"

$ws.Range("A14").Value = "process_model_id"
$ws.Range("B14").Value = "str"

$ws.Range("A16").Value = "SubProcessID"
$ws.Range("B16").Value = "str"

$ws.Range("A17").Value = "time:timestamp"
$ws.Range("B17").Value = "datetime"

$ws.Range("A18").Value = "requested_service_url"
$ws.Range("B18").Value = "str"

$ws.Range("A19").Value = "lifecycle:transition"
$ws.Range("B19").Value = "str"

$ws.Range("A20").Value = "case"
$ws.Range("B20").Value = "str"
